$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (this is a new daily price entry),
# shifting existing rows 85-175 down to 86-176.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new data record.
$ws.Cells.Item(85, 1).Value = 4
$ws.Cells.Item(85, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(85, 3).Value = "Los Lagos"
$ws.Cells.Item(85, 4).Value = 44601
$ws.Cells.Item(85, 5).Value = 10
$ws.Cells.Item(85, 6).Value = 100112039
$ws.Cells.Item(85, 7).Value = "Ciboulette"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 40
$ws.Cells.Item(85, 11).Value = 2500
$ws.Cells.Item(85, 12).Value = 2500
$ws.Cells.Item(85, 13).Value = 2500
$ws.Cells.Item(85, 14).Value = "`$/docena de atados"
$ws.Cells.Item(85, 15).Value = "Región Metropolitana"
$ws.Cells.Item(85, 16).Value = 833
$ws.Cells.Item(85, 17).Value = 3
$ws.Cells.Item(85, 18).Value = "Hortaliza"
